# Generate Report for handoff
#
# Semantic edit being applied (per commit "Generate Report for handoff"):
#   1. The source file 3cef2cad-cc9e-4017-ad68-c3673c44d788.md was renamed to
#      fe439221-6ba3-4a48-835d-04cb71d58f08.md; its handoff package / timestamp
#      for both the zh-cn and de-de targets was refreshed accordingly.
#   2. A brand-new source file 869b81ce-d447-418b-8544-3ce92f82222c.md was
#      picked up whose handoff transform failed, so a new row is inserted for
#      it (ahead of the static ".localization-config" row, which shifts down
#      one row on every sheet).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "fe439221-6ba3-4a48-835d-04cb71d58f08.md"

$ws1.Rows.Item(3).Insert()
$ws1.Range("A3").Value = "869b81ce-d447-418b-8544-3ce92f82222c.md"
$ws1.Range("B3").Value = "Handoff transform failed"
$ws1.Range("C3").Value = "Handoff transform failed"

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/42edd7ac7443c0f5124194abe45f2049552ac8e4/e2e/fe439221-6ba3-4a48-835d-04cb71d58f08.md", "", "", "fe439221-6ba3-4a48-835d-04cb71d58f08.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/42edd7ac7443c0f5124194abe45f2049552ac8e4/e2e/869b81ce-d447-418b-8544-3ce92f82222c.md", "", "", "869b81ce-d447-418b-8544-3ce92f82222c.md")
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/42edd7ac7443c0f5124194abe45f2049552ac8e4/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "fe439221-6ba3-4a48-835d-04cb71d58f08.md"
$ws2.Range("C2").Value = "fe439221-6ba3-4a48-835d-04cb71d58f08.ca67c59f4598e26446334c97a7b40b376a7b7ae5.zh-cn.xlf"
$ws2.Range("D2").Value = "2016-01-26 12:10:51"

$ws2.Rows.Item(3).Insert()
$ws2.Range("C3").Clear()
$ws2.Range("A3").Value = "869b81ce-d447-418b-8544-3ce92f82222c.md"
$ws2.Range("B3").Value = "Handoff transform failed"
$ws2.Range("D3").Value = "0001-01-01 00:00:00"
$ws2.Range("G3").Value = "0001-01-01 00:00:00"
$ws2.Range("H3").Value = "Ignored"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/42edd7ac7443c0f5124194abe45f2049552ac8e4/e2e/fe439221-6ba3-4a48-835d-04cb71d58f08.md", "", "", "fe439221-6ba3-4a48-835d-04cb71d58f08.md")
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6ccfc053f2feb09a2b1b1eb47af446ae6706576e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/fe439221-6ba3-4a48-835d-04cb71d58f08.ca67c59f4598e26446334c97a7b40b376a7b7ae5.zh-cn.xlf", "", "", "fe439221-6ba3-4a48-835d-04cb71d58f08.ca67c59f4598e26446334c97a7b40b376a7b7ae5.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/42edd7ac7443c0f5124194abe45f2049552ac8e4/e2e/869b81ce-d447-418b-8544-3ce92f82222c.md", "", "", "869b81ce-d447-418b-8544-3ce92f82222c.md")
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/42edd7ac7443c0f5124194abe45f2049552ac8e4/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "fe439221-6ba3-4a48-835d-04cb71d58f08.md"
$ws3.Range("C2").Value = "fe439221-6ba3-4a48-835d-04cb71d58f08.ca67c59f4598e26446334c97a7b40b376a7b7ae5.de-de.xlf"
$ws3.Range("D2").Value = "2016-01-26 12:11:02"

$ws3.Rows.Item(3).Insert()
$ws3.Range("C3").Clear()
$ws3.Range("A3").Value = "869b81ce-d447-418b-8544-3ce92f82222c.md"
$ws3.Range("B3").Value = "Handoff transform failed"
$ws3.Range("D3").Value = "0001-01-01 00:00:00"
$ws3.Range("G3").Value = "0001-01-01 00:00:00"
$ws3.Range("H3").Value = "Ignored"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/42edd7ac7443c0f5124194abe45f2049552ac8e4/e2e/fe439221-6ba3-4a48-835d-04cb71d58f08.md", "", "", "fe439221-6ba3-4a48-835d-04cb71d58f08.md")
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d0724ee4709583fdbce5de62de5137f196c38d5b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/fe439221-6ba3-4a48-835d-04cb71d58f08.ca67c59f4598e26446334c97a7b40b376a7b7ae5.de-de.xlf", "", "", "fe439221-6ba3-4a48-835d-04cb71d58f08.ca67c59f4598e26446334c97a7b40b376a7b7ae5.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/42edd7ac7443c0f5124194abe45f2049552ac8e4/e2e/869b81ce-d447-418b-8544-3ce92f82222c.md", "", "", "869b81ce-d447-418b-8544-3ce92f82222c.md")
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/42edd7ac7443c0f5124194abe45f2049552ac8e4/.localization-config", "", "", ".localization-config")
